# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# For each data row (16-55) of Hoja1:
#   - Column E ("Periodo Mora") gets the period label, now running in
#     ascending chronological order (1612 .. 2003) instead of descending.
#   - Column F ("Valor Mora") / Column G ("Salario Basico") get refreshed
#     amounts.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")
if (-not $ws) { $ws = $wb.ActiveSheet }

$rows = @(
  @{Row=16; Period="1612"; F=27578; G=781242},
  @{Row=17; Period="1701"; F=27578; G=781242},
  @{Row=18; Period="1702"; F=27578; G=781242},
  @{Row=19; Period="1703"; F=27578; G=781242},
  @{Row=20; Period="1704"; F=27578; G=781242},
  @{Row=21; Period="1705"; F=27578; G=781242},
  @{Row=22; Period="1706"; F=27578; G=781242},
  @{Row=23; Period="1707"; F=27578; G=781242},
  @{Row=24; Period="1708"; F=27578; G=781242},
  @{Row=25; Period="1709"; F=27578; G=781242},
  @{Row=26; Period="1710"; F=27578; G=781242},
  @{Row=27; Period="1711"; F=27578; G=781242},
  @{Row=28; Period="1712"; F=27578; G=781242},
  @{Row=29; Period="1801"; F=27578; G=781242},
  @{Row=30; Period="1802"; F=27578; G=781242},
  @{Row=31; Period="1803"; F=27578; G=781242},
  @{Row=32; Period="1804"; F=27578; G=781242},
  @{Row=33; Period="1805"; F=27578; G=781242},
  @{Row=34; Period="1806"; F=27578; G=781242},
  @{Row=35; Period="1807"; F=27578; G=781242},
  @{Row=36; Period="1808"; F=27578; G=781242},
  @{Row=37; Period="1809"; F=31249; G=781242},
  @{Row=38; Period="1810"; F=31249; G=781242},
  @{Row=39; Period="1811"; F=31249; G=781242},
  @{Row=40; Period="1812"; F=31249; G=781242},
  @{Row=41; Period="1901"; F=31249; G=781242},
  @{Row=42; Period="1902"; F=31249; G=781242},
  @{Row=43; Period="1903"; F=31249; G=781242},
  @{Row=44; Period="1904"; F=31249; G=781242},
  @{Row=45; Period="1905"; F=31249; G=781242},
  @{Row=46; Period="1906"; F=31249; G=781242},
  @{Row=47; Period="1907"; F=31249; G=781242},
  @{Row=48; Period="1908"; F=31249; G=781242},
  @{Row=49; Period="1909"; F=31249; G=781242},
  @{Row=50; Period="1910"; F=31249; G=781242},
  @{Row=51; Period="1911"; F=31249; G=781242},
  @{Row=52; Period="1912"; F=31249; G=781242},
  @{Row=53; Period="2001"; F=31249; G=781242},
  @{Row=54; Period="2002"; F=31249; G=781242},
  @{Row=55; Period="2003"; F=31249; G=781242}
)

foreach ($r in $rows) {
    $ws.Range("E" + $r.Row).Value = $r.Period
    $ws.Range("F" + $r.Row).Value = $r.F
    $ws.Range("G" + $r.Row).Value = $r.G
}
